# Adjusts the "autodiagnostico" test-case sheet:
#  - renames the header in column B from "Nombre Caso" to "Nombre/ Descripcion"
#  - widens column G so the longer "Pasos a Seguir" text is easier to read
#  - moves the active selection to the first data row (B2:M2) as left by the editor

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the column B header text
$headerCell = $ws.Range("B1")
$headerCell.Value = "Nombre/ Descripcion"

# Re-apply the header-style border around B1 (thin box, no bottom rule) without
# introducing any stray/duplicate border definitions: BorderAround first lays
# down a uniform thin black box, then we strip the bottom edge to match the
# rest of the header row's "open bottom" styling.
$headerCell.BorderAround(1, 2, 1, 0)
$headerCell.Borders.Item(9).LineStyle = -4142

# 2. Widen column G (Pasos a Seguir) by 15 characters
$colG = $ws.Columns.Item(7)
$colG.ColumnWidth = $colG.ColumnWidth + 15

# 3. Leave the selection on the first data row, matching the editor's last action
$ws.Range("B2:M2").Select()
